$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "NA" under duplicate_image_filename (column E) for data rows 2-21
$ws.Range("E2:E21").Value = "NA"
